# Update "想去人数" (F column) figures and one cover image URL (I29)
# on both the "展览" and "全部类型" worksheets, which contain
# identical data tables.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$fUpdates = @{
    2  = 1069
    3  = 764
    6  = 1101
    8  = 1804
    9  = 6752
    13 = 109
    14 = 382
    15 = 140
    16 = 6954
    18 = 1303
    19 = 140
    20 = 119
    22 = 113
    30 = 598
    32 = 81
    34 = 62
    36 = 64
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }

    $ws.Cells.Item(29, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/talOodLW1714030986517.png"
}
